$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(339).Insert()

$ws.Range("A339").Value = 10
$ws.Range("B339").Value = "Vega Modelo de Temuco"
$ws.Range("C339").Value = "La Araucanía"
$ws.Range("D339").Value = 44855
$ws.Range("E339").Value = 9
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100108
$ws.Range("H339").Value = "Tropicales y subtropicales"
$ws.Range("I339").Value = 100108002
$ws.Range("J339").Value = "Mango"
$ws.Range("K339").Value = "Sin especificar"
$ws.Range("L339").Value = "Primera"
$ws.Range("M339").Value = 430
$ws.Range("N339").Value = 8000
$ws.Range("O339").Value = 9000
$ws.Range("P339").Value = 8419
$ws.Range("Q339").Value = "`$/bandeja 4 kilos"
$ws.Range("R339").Value = "Brasil"
$ws.Range("S339").Value = 2105
$ws.Range("T339").Value = 4

$ws.Range("D339").NumberFormat = $ws.Range("D340").NumberFormat
